$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 0.04741066666666666
$ws.Range("H2").Value = 0.142232
$ws.Range("I2").Value = 0.003188134523263584
$ws.Range("J2").Value = 0.003188134523263585
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 3.197736333333333
$ws.Range("N2").Value = 9.593208999999998
$ws.Range("Q2").Value = 0.1516068113875555
$ws.Range("R2").Value = 1.364461302488
$ws.Range("S2").Value = 0.003188134523263584
$ws.Range("T2").Value = 0.003188134523263585

$ws.Range("I3").Value = 0.01595759596384214
$ws.Range("J3").Value = 0.01595759596384214
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 3.197736333333333
$ws.Range("N3").Value = 9.593208999999998
$ws.Range("Q3").Value = 0.7588388205816665
$ws.Range("R3").Value = 6.829549385234999
$ws.Range("S3").Value = 0.01595759596384214
$ws.Range("T3").Value = 0.01595759596384214

$ws.Range("I4").Value = 0.9808542695128942
$ws.Range("J4").Value = 0.9808542695128943
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 3.197736333333333
$ws.Range("N4").Value = 9.593208999999998
$ws.Range("Q4").Value = 46.64300930579821
$ws.Range("R4").Value = 419.7870837521839
$ws.Range("S4").Value = 0.9808542695128942
$ws.Range("T4").Value = 0.9808542695128943
